{"js": "// Replace the date line and the two-digit multiplication problems per the\n// commit's diff. Every old value in the document is unique, so a simple\n// search-and-replace keyed on the exact old text is unambiguous.\nconst replacements = [\n  [\"2025-07-16 Wednesday\", \"2025-07-17 Thursday\"],\n  [\"70\u00d755=\", \"72\u00d756=\"],\n  [\"91\u00d736=\", \"54\u00d731=\"],\n  [\"92\u00d715=\", \"39\u00d734=\"],\n  [\"15\u00d717=\", \"41\u00d788=\"],\n  [\"54\u00d792=\", \"60\u00d753=\"],\n  [\"21\u00d749=\", \"98\u00d713=\"],\n  [\"87\u00d768=\", \"82\u00d729=\"],\n  [\"79\u00d720=\", \"91\u00d791=\"],\n  [\"31\u00d719=\", \"80\u00d715=\"],\n  [\"60\u00d794=\", \"65\u00d769=\"],\n  [\"69\u00d753=\", \"35\u00d797=\"],\n  [\"90\u00d725=\", \"15\u00d741=\"],\n  [\"56\u00d784=\", \"54\u00d730=\"],\n  [\"99\u00d719=\", \"66\u00d733=\"],\n  [\"11\u00d789=\", \"73\u00d733=\"],\n  [\"26\u00d753=\", \"94\u00d769=\"],\n  [\"92\u00d772=\", \"39\u00d770=\"],\n  [\"22\u00d721=\", \"16\u00d749=\"],\n  [\"17\u00d759=\", \"49\u00d769=\"],\n  [\"22\u00d743=\", \"86\u00d734=\"],\n  [\"47\u00d715=\", \"77\u00d755=\"],\n  [\"56\u00d778=\", \"93\u00d737=\"],\n  [\"48\u00d761=\", \"53\u00d768=\"],\n  [\"91\u00d772=\", \"13\u00d746=\"],\n  [\"27\u00d721=\", \"31\u00d790=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the two-digit multiplication problems per the\n# commit's diff. Every old value in the document is unique, so Find/Replace\n# keyed on the exact old text (MatchCase, whole-document scope) is\n# unambiguous and safe to run as a straight ReplaceAll for each pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-07-16 Wednesday\", \"2025-07-17 Thursday\"),\n    @(\"70\u00d755=\", \"72\u00d756=\"),\n    @(\"91\u00d736=\", \"54\u00d731=\"),\n    @(\"92\u00d715=\", \"39\u00d734=\"),\n    @(\"15\u00d717=\", \"41\u00d788=\"),\n    @(\"54\u00d792=\", \"60\u00d753=\"),\n    @(\"21\u00d749=\", \"98\u00d713=\"),\n    @(\"87\u00d768=\", \"82\u00d729=\"),\n    @(\"79\u00d720=\", \"91\u00d791=\"),\n    @(\"31\u00d719=\", \"80\u00d715=\"),\n    @(\"60\u00d794=\", \"65\u00d769=\"),\n    @(\"69\u00d753=\", \"35\u00d797=\"),\n    @(\"90\u00d725=\", \"15\u00d741=\"),\n    @(\"56\u00d784=\", \"54\u00d730=\"),\n    @(\"99\u00d719=\", \"66\u00d733=\"),\n    @(\"11\u00d789=\", \"73\u00d733=\"),\n    @(\"26\u00d753=\", \"94\u00d769=\"),\n    @(\"92\u00d772=\", \"39\u00d770=\"),\n    @(\"22\u00d721=\", \"16\u00d749=\"),\n    @(\"17\u00d759=\", \"49\u00d769=\"),\n    @(\"22\u00d743=\", \"86\u00d734=\"),\n    @(\"47\u00d715=\", \"77\u00d755=\"),\n    @(\"56\u00d778=\", \"93\u00d737=\"),\n    @(\"48\u00d761=\", \"53\u00d768=\"),\n    @(\"91\u00d772=\", \"13\u00d746=\"),\n    @(\"27\u00d721=\", \"31\u00d790=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$newText, [ref]\"wdReplaceAll\")\n}\n\n$d.Save()\n"}
